$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 209.5
$ws.Range("J42").Value = 49.666668
$ws.Range("L42").Value = 149.000004
$ws.Range("N42").Value = -609.000004
$ws.Range("H58").Value = 378
$ws.Range("I58").Value = 111.42857
$ws.Range("K58").Value = 334.28571
$ws.Range("M58").Value = -184.28571
$ws.Range("H70").Value = 2473.4736
$ws.Range("I70").Value = 2131.3333
$ws.Range("J70").Value = 2537.625
$ws.Range("K70").Value = 6393.999899999999
$ws.Range("L70").Value = 7612.875
$ws.Range("M70").Value = -6123.999899999999
$ws.Range("N70").Value = -8152.875
$ws.Range("H73").Value = 2473.4736
$ws.Range("I73").Value = 2131.3333
$ws.Range("J73").Value = 2537.625
$ws.Range("K73").Value = 6393.999899999999
$ws.Range("L73").Value = 7612.875
$ws.Range("M73").Value = -5457.999899999999
$ws.Range("N73").Value = -9484.875
$ws.Range("H98").Value = 1141.5714
$ws.Range("I98").Value = 1141.5714
$ws.Range("K98").Value = 1141.5714
$ws.Range("M98").Value = 356.4286
$ws.Range("H103").Value = 358.08334
$ws.Range("J103").Value = 270
$ws.Range("L103").Value = 810
$ws.Range("N103").Value = -1982
$ws.Range("H115").Value = 184.66667
$ws.Range("I115").Value = 184.66667
$ws.Range("K115").Value = 554.00001
$ws.Range("M115").Value = 1012.99999
$ws.Range("H122").Value = 1141.5714
$ws.Range("I122").Value = 1141.5714
$ws.Range("K122").Value = 3424.7142
$ws.Range("M122").Value = -974.7142000000003
$ws.Range("H132").Value = 1755.1538
$ws.Range("I132").Value = 1417.875
$ws.Range("K132").Value = 4253.625
$ws.Range("M132").Value = -1723.625
$ws.Range("H137").Value = 2811.7778
$ws.Range("I137").Value = 1901.3
$ws.Range("K137").Value = 5703.9
$ws.Range("M137").Value = -3153.9
$ws.Range("H138").Value = 15386460
$ws.Range("I138").Value = 965.63635
$ws.Range("J138").Value = 47622732
$ws.Range("K138").Value = 2896.90905
$ws.Range("L138").Value = 142868196
$ws.Range("M138").Value = 2243.09095
$ws.Range("N138").Value = -142878476
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12816.848
$ws.Range("I32").Value = 3688.8704
$ws.Range("K32").Value = 3688.8704
$ws.Range("M32").Value = -3401.8704
$ws.Range("H50").Value = 1787.5
$ws.Range("J50").Value = 2116.6667
$ws.Range("L50").Value = 2116.6667
$ws.Range("N50").Value = -3544.6667
$ws.Range("H54").Value = 10400
$ws.Range("I54").Value = 8000
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 8000
$ws.Range("L54").Value = 20000
$ws.Range("M54").Value = -7231
$ws.Range("N54").Value = -21538
$ws.Range("H61").Value = 20011146
$ws.Range("I61").Value = 26326244
$ws.Range("K61").Value = 26326244
$ws.Range("M61").Value = -26326032
$ws.Range("H102").Value = 2473.4062
$ws.Range("I102").Value = 1937.3914
$ws.Range("J102").Value = 3843.2222
$ws.Range("K102").Value = 1937.3914
$ws.Range("L102").Value = 3843.2222
$ws.Range("M102").Value = -315.3914
$ws.Range("N102").Value = -7087.2222
$ws.Range("H122").Value = 1572.25
$ws.Range("I122").Value = 1263.1666
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 3789.4998
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -1339.4998
$ws.Range("N122").Value = -12398.5
$ws.Range("H132").Value = 2557.3809
$ws.Range("I132").Value = 2567.054
$ws.Range("J132").Value = 2485.8
$ws.Range("K132").Value = 7701.162
$ws.Range("L132").Value = 7457.400000000001
$ws.Range("M132").Value = -5171.162
$ws.Range("N132").Value = -12517.4
$ws.Range("H136").Value = 20011146
$ws.Range("I136").Value = 26326244
$ws.Range("K136").Value = 78978732
$ws.Range("M136").Value = -78976182
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 62820
$ws.Range("J2").Value = 69230
$ws.Range("L2").Value = 69230
$ws.Range("N2").Value = -69456
$ws.Range("H92").Value = 49999.5
$ws.Range("J92").Value = 49999.5
$ws.Range("L92").Value = 49999.5
$ws.Range("N92").Value = -54991.5
$ws.Range("H134").Value = 4036.9302
$ws.Range("I134").Value = 4106.5674
$ws.Range("K134").Value = 12319.7022
$ws.Range("M134").Value = -9784.7022
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 49999.668
$ws.Range("J9").Value = 49999.668
$ws.Range("L9").Value = 49999.668
$ws.Range("N9").Value = -50335.668
$ws.Range("H31").Value = 3771.7646
$ws.Range("I31").Value = 2123.5652
$ws.Range("K31").Value = 2123.5652
$ws.Range("M31").Value = -1828.5652
$ws.Range("H34").Value = 3771.7646
$ws.Range("I34").Value = 2123.5652
$ws.Range("K34").Value = 2123.5652
$ws.Range("M34").Value = -1921.5652
$ws.Range("H68").Value = 41406.668
$ws.Range("J68").Value = 41406.668
$ws.Range("L68").Value = 41406.668
$ws.Range("N68").Value = -42904.668
$ws.Range("H71").Value = 41406.668
$ws.Range("J71").Value = 41406.668
$ws.Range("L71").Value = 124220.004
$ws.Range("N71").Value = -131708.004
$ws.Range("H105").Value = 1820.5834
$ws.Range("I105").Value = 1394.8
$ws.Range("K105").Value = 1394.8
$ws.Range("M105").Value = 352.2
$ws.Range("H107").Value = 1008.73914
$ws.Range("I107").Value = 1028.579
$ws.Range("J107").Value = 914.5
$ws.Range("K107").Value = 1028.579
$ws.Range("L107").Value = 914.5
$ws.Range("M107").Value = 891.421
$ws.Range("N107").Value = -4754.5
$ws.Range("H132").Value = 2637.3428
$ws.Range("I132").Value = 2844.3103
$ws.Range("K132").Value = 8532.930899999999
$ws.Range("M132").Value = -6002.930899999999
$ws.Range("H134").Value = 3320.6
$ws.Range("I134").Value = 2478.4
$ws.Range("K134").Value = 7435.200000000001
$ws.Range("M134").Value = -4900.200000000001
$ws.Range("H141").Value = 36800
$ws.Range("J141").Value = 36800
$ws.Range("L141").Value = 36800
$ws.Range("N141").Value = -47160
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1035.625
$ws.Range("J113").Value = 1419.75
$ws.Range("L113").Value = 4259.25
$ws.Range("N113").Value = -8599.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5524.875
$ws.Range("J102").Value = 12000
$ws.Range("L102").Value = 12000
$ws.Range("N102").Value = -15244
$ws.Range("H113").Value = 366810.72
$ws.Range("I113").Value = 573585
$ws.Range("J113").Value = 4955.75
$ws.Range("K113").Value = 573585
$ws.Range("L113").Value = 4955.75
$ws.Range("M113").Value = -571415
$ws.Range("N113").Value = -9295.75
$ws.Range("H122").Value = 1887.4
$ws.Range("I122").Value = 1942.6
$ws.Range("K122").Value = 5827.799999999999
$ws.Range("M122").Value = -3377.799999999999
$ws.Range("H132").Value = 1990.7142
$ws.Range("I132").Value = 1916.2963
$ws.Range("K132").Value = 5748.8889
$ws.Range("M132").Value = -3218.8889
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3424
$ws.Range("I7").Value = 3424
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3424
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = -3312
$ws.Range("H40").Value = 3497.1667
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = $null
$ws.Range("H122").Value = 3189.5557
$ws.Range("I122").Value = 2959.4285
$ws.Range("J122").Value = 3995
$ws.Range("K122").Value = 8878.2855
$ws.Range("L122").Value = 11985
$ws.Range("M122").Value = -6428.2855
$ws.Range("N122").Value = -16885
$ws.Range("H126").Value = 3424
$ws.Range("I126").Value = 3424
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10272
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = -7802
$ws.Range("H132").Value = 7808.3438
$ws.Range("I132").Value = 7569.9453
$ws.Range("K132").Value = 22709.8359
$ws.Range("M132").Value = -20179.8359
$ws.Range("H136").Value = 3830.8958
$ws.Range("J136").Value = 2613
$ws.Range("L136").Value = 7839
$ws.Range("N136").Value = -12939
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7482.3335
$ws.Range("I15").Value = 7481
$ws.Range("J15").Value = 7483.6665
$ws.Range("K15").Value = 7481
$ws.Range("L15").Value = 7483.6665
$ws.Range("M15").Value = -7193
$ws.Range("N15").Value = -8059.6665
$ws.Range("H81").Value = 1662.9412
$ws.Range("I81").Value = 848
$ws.Range("K81").Value = 1696
$ws.Range("M81").Value = -635
$ws.Range("H84").Value = 1662.9412
$ws.Range("I84").Value = 848
$ws.Range("K84").Value = 8480
$ws.Range("M84").Value = -3176
$ws.Range("H122").Value = 3039.7778
$ws.Range("I122").Value = 1622.381
$ws.Range("K122").Value = 4867.143
$ws.Range("M122").Value = -2417.143
$ws.Range("H132").Value = 2201.7068
$ws.Range("I132").Value = 2324.2449
$ws.Range("J132").Value = 1534.5555
$ws.Range("K132").Value = 6972.734700000001
$ws.Range("L132").Value = 4603.666499999999
$ws.Range("M132").Value = -4442.734700000001
$ws.Range("N132").Value = -9663.666499999999
$ws.Range("H136").Value = 5660.8335
$ws.Range("I136").Value = 3679.1875
$ws.Range("K136").Value = 11037.5625
$ws.Range("M136").Value = -8487.5625
